$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column C and D
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 1
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 2
$ws.Range("C6").Value = 7
$ws.Range("D6").Value = 7
$ws.Range("D7").Value = 8

# Update the active selection to C7
$ws.Range("C7").Select()
